# "Daily 100 Error Counts" — drop the last two data rows (25 & 26), which
# get cleared back out to the same empty/date-formatted shape row 27
# already had, and update the saved view state (scroll position + active
# selection) to where the author last left the cursor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 (45971 / 598 / 20 / 578) and row 26 (45973 / 680 / 26 / 654) are
# cleared back to an empty cell — ClearContents() wipes the values but
# keeps the date number-format style already applied to column A.
$ws.Range("A25:D25").ClearContents() | Out-Null
$ws.Range("A26:D26").ClearContents() | Out-Null

# Restore the view: scrolled so row 8 is at the top, with G21 selected.
$ws.Activate() | Out-Null
$ws.Range("A8").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("G21").Select() | Out-Null
